$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = "[파이썬 독학] 파이썬 생활 밀착형 프로젝트 X 클래스101 온라인 강의 오픈 (w/쿠폰 3만원)"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/614"

$ws.Range("D24").Value = "[근황] 졸업, 억셉, 첫 출근"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222662292675"

$ws.Range("D29").Value = "[만화] 인턴일기 72~80"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-11/"
